# Insert a new record as row 14 ("Hortaliza, Terminal Hortofrutícola Agro
# Chillán - Cebollín" weekly update). This shifts the existing rows 14-110
# down to 15-111 (Excel updates all the relative references automatically)
# and grows the sheet from A1:R110 to A1:R111.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(14).Insert()

$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C14").Value = "Ñuble"
$ws.Range("D14").Value = 44991
$ws.Range("E14").Value = 16
$ws.Range("F14").Value = 100112037
$ws.Range("G14").Value = "Cebollín"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 6000
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = 6000
$ws.Range("N14").Value = "`$/paquete 36 unidades"
$ws.Range("O14").Value = "Provincia de Diguillín"
$ws.Range("P14").Value = 167
$ws.Range("Q14").Value = 36
$ws.Range("R14").Value = "Hortaliza"
